$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.951.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.414.52'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '410.20'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.46'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.635'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.57%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +10.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.66'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.07%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.09'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.03%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.141'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.20'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.14%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.937.21'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000206'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +44.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.460.68'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.16'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.08'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +6.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '61.868.22'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '456.79'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +46.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.37'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +9.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.16'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.25'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.66'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +13.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +8.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.77'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.74'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.57'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.97'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.66%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '42.94'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.168'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0497'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.28'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.997'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.134'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.314'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.85%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.29'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.23%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '141.11'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +7.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.49'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.62'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +6.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.746.71'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.111.50'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.53%  '
